$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new weekly record is inserted at row 134, pushing the existing
# records (previously rows 134-165) down by one row (to rows 135-166).
$ws.Rows.Item(134).Insert()

# Populate the newly inserted row 134. Its contents mirror what used to be
# row 134 (now shifted to row 135), except for a new date (one week later).
$ws.Range("A134").Value = 5
$ws.Range("B134").Value = "Macroferia Regional de Talca"
$ws.Range("C134").Value = "Maule"
$ws.Range("D134").Value = 44511
$ws.Range("E134").Value = 7
$ws.Range("F134").Value = 100112008
$ws.Range("G134").Value = "Coliflor"
$ws.Range("H134").Value = "Sin especificar"
$ws.Range("I134").Value = "Primera"
$ws.Range("J134").Value = 4000
$ws.Range("K134").Value = 600
$ws.Range("L134").Value = 600
$ws.Range("M134").Value = 600
$ws.Range("N134").Value = "$/unidad"
$ws.Range("O134").Value = "Región del Maule"
$ws.Range("P134").Value = 600
$ws.Range("Q134").Value = 1
$ws.Range("R134").Value = "Hortaliza"
